$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text (11/12/2016 ->
#    11/16/2016). This automatically-updating field is rendered from the
#    Date Placeholder that lives on the slide master AND on every individual
#    slide layout, so every one of those copies has to be touched.
# ---------------------------------------------------------------------------
$newDate = "11/16/2016"
$ppPlaceholderDate = 16

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    $phType = -1
    try { $phType = $shp.PlaceholderFormat.Type } catch { $phType = -1 }
    if ($phType -eq $ppPlaceholderDate) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $shp = $cl.Shapes.Item($i)
        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch { $phType = -1 }
        if ($phType -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Rename the "Gap" label to "IPG" on slide 1 (inside the Ethernet frame
#    diagram group) - IPG is the common abbreviation for inter-packet gap.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Gap") {
            $shp.TextFrame.TextRange.Text = "IPG"
        }
    }
}
